$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.161.06"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.394.22"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.37"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.88"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.51"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "2.755.80"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "2.382.66"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "43.199.24"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.03"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "0.0₃0897"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.28"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.91"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.65"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.20"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.11"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  +13.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.24"
$ws.Range("E34").Value = "  +5.72%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "130.92"
$ws.Range("E36").Value = "  +16.33%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.90"
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.08"
$ws.Range("E42").Value = "  -6.81%  "
$ws.Range("D43").Value = "1.941.92"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.14"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("E47").Value = "  -7.33%  "
$ws.Range("D48").Value = "2.609.33"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.26"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.93"
$ws.Range("E51").Value = "  -0.31%  "
